$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Máster en Neuropsicología Clínica / 2023-2024 / Leidy Nathaly Peláez Bernal / VIU link
$ws.Range("A2").Value = "Máster en Neuropsicología Clínica"
$ws.Range("B2").Value = "2023-2024"
$ws.Range("C2").Value = "Leidy Nathaly Peláez Bernal"
$ws.Range("D2").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"

# Row 3: Máster en Neuropsicología Clínica / 2023-2024 / Jimena Zanizo Chambi / VIU link
$ws.Range("A3").Value = "Máster en Neuropsicología Clínica"
$ws.Range("B3").Value = "2023-2024"
$ws.Range("C3").Value = "Jimena Zanizo Chambi"
$ws.Range("D3").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"

# Row 4: Máster en Neuropsicología Clínica / 2023-2024 / Liceth Andrea Zaraza Osorio / VIU link
$ws.Range("A4").Value = "Máster en Neuropsicología Clínica"
$ws.Range("B4").Value = "2023-2024"
$ws.Range("C4").Value = "Liceth Andrea Zaraza Osorio"
$ws.Range("D4").Value = "\href{https://www.universidadviu.com/co/}{Universidad Internacional de Valencia}, España"

# Row 4 now holds wrapped text like rows 2-3, so match their auto-fitted height
$ws.Rows.Item(4).RowHeight = 43.2

# Update the selection to match the published state
$ws.Range("A2:D4").Select()
